# Update Betfair Back/Lay odds sheet for 2026-01-05
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Benfica B x Porto B
$ws.Range("H3").Value = 4
$ws.Range("Q3").Value = 1.69

# Row 4 - Maccabi Haifa x Hapoel Haifa
$ws.Range("G4").Value = 1.48

# Row 5 - now Triestina x ASD Alcione
$ws.Range("D5").Value = "Triestina"
$ws.Range("E5").Value = "ASD Alcione"
$ws.Range("F5").Value = 2.08
$ws.Range("G5").Value = 2.34
$ws.Range("H5").Value = 3.55
$ws.Range("I5").Value = 4.4
$ws.Range("J5").Value = 3.2
$ws.Range("K5").Value = 3.75
$ws.Range("P5").Value = 1.7
$ws.Range("Q5").Value = 2.1

# Row 6 - now Cosenza x Monopoli
$ws.Range("D6").Value = "Cosenza"
$ws.Range("E6").Value = "Monopoli"
$ws.Range("F6").Value = 1.81
$ws.Range("G6").Value = 2.08
$ws.Range("H6").Value = 4.9
$ws.Range("I6").Value = 6.6
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 3.55
$ws.Range("P6").Value = 1.59
$ws.Range("Q6").Value = 2.16

# Row 7 - now Benevento x Crotone
$ws.Range("D7").Value = "Benevento"
$ws.Range("E7").Value = "Crotone"
$ws.Range("F7").Value = 1.7
$ws.Range("G7").Value = 1.83
$ws.Range("H7").Value = 4.5
$ws.Range("I7").Value = 7.2
$ws.Range("J7").Value = 3.7
$ws.Range("K7").Value = 4.8
$ws.Range("P7").Value = 1.94
$ws.Range("Q7").Value = 1.89

# Row 8 - Montpellier x Dunkerque
$ws.Range("F8").Value = 2.32
$ws.Range("G8").Value = 2.6
$ws.Range("H8").Value = 3.15
$ws.Range("I8").Value = 3.6
$ws.Range("J8").Value = 3.15
$ws.Range("K8").Value = 3.6

# Row 9 - Vizela x Torreense
$ws.Range("F9").Value = 1.89
$ws.Range("G9").Value = 2.12
$ws.Range("H9").Value = 4.2
$ws.Range("I9").Value = 5.7
$ws.Range("K9").Value = 3.85
